$d = $word.ActiveDocument

# Phase 1: replace each original expression with a unique placeholder token
# to avoid collisions where a new value equals another cell's original value
# (e.g. "4+48=" -> "25+70=" while another cell's original "25+70=" -> "93-59=").
$d.Content.Find.Execute("70+19=", $true, $false, $false, $false, $false, $true, 1, $false, "@@0@@", 2) | Out-Null
$d.Content.Find.Execute("24+74=", $true, $false, $false, $false, $false, $true, 1, $false, "@@1@@", 2) | Out-Null
$d.Content.Find.Execute("4+47=", $true, $false, $false, $false, $false, $true, 1, $false, "@@2@@", 2) | Out-Null
$d.Content.Find.Execute("32+22=", $true, $false, $false, $false, $false, $true, 1, $false, "@@3@@", 2) | Out-Null
$d.Content.Find.Execute("51-0=", $true, $false, $false, $false, $false, $true, 1, $false, "@@4@@", 2) | Out-Null
$d.Content.Find.Execute("83-32=", $true, $false, $false, $false, $false, $true, 1, $false, "@@5@@", 2) | Out-Null
$d.Content.Find.Execute("13+42=", $true, $false, $false, $false, $false, $true, 1, $false, "@@6@@", 2) | Out-Null
$d.Content.Find.Execute("26+20=", $true, $false, $false, $false, $false, $true, 1, $false, "@@7@@", 2) | Out-Null
$d.Content.Find.Execute("85-28=", $true, $false, $false, $false, $false, $true, 1, $false, "@@8@@", 2) | Out-Null
$d.Content.Find.Execute("10-8=", $true, $false, $false, $false, $false, $true, 1, $false, "@@9@@", 2) | Out-Null
$d.Content.Find.Execute("47-23=", $true, $false, $false, $false, $false, $true, 1, $false, "@@10@@", 2) | Out-Null
$d.Content.Find.Execute("49+50=", $true, $false, $false, $false, $false, $true, 1, $false, "@@11@@", 2) | Out-Null
$d.Content.Find.Execute("22+69=", $true, $false, $false, $false, $false, $true, 1, $false, "@@12@@", 2) | Out-Null
$d.Content.Find.Execute("60+23=", $true, $false, $false, $false, $false, $true, 1, $false, "@@13@@", 2) | Out-Null
$d.Content.Find.Execute("99-70=", $true, $false, $false, $false, $false, $true, 1, $false, "@@14@@", 2) | Out-Null
$d.Content.Find.Execute("0+34=", $true, $false, $false, $false, $false, $true, 1, $false, "@@15@@", 2) | Out-Null
$d.Content.Find.Execute("58-41=", $true, $false, $false, $false, $false, $true, 1, $false, "@@16@@", 2) | Out-Null
$d.Content.Find.Execute("15+20=", $true, $false, $false, $false, $false, $true, 1, $false, "@@17@@", 2) | Out-Null
$d.Content.Find.Execute("26+2=", $true, $false, $false, $false, $false, $true, 1, $false, "@@18@@", 2) | Out-Null
$d.Content.Find.Execute("62-28=", $true, $false, $false, $false, $false, $true, 1, $false, "@@19@@", 2) | Out-Null
$d.Content.Find.Execute("64-38=", $true, $false, $false, $false, $false, $true, 1, $false, "@@20@@", 2) | Out-Null
$d.Content.Find.Execute("54+2=", $true, $false, $false, $false, $false, $true, 1, $false, "@@21@@", 2) | Out-Null
$d.Content.Find.Execute("39-34=", $true, $false, $false, $false, $false, $true, 1, $false, "@@22@@", 2) | Out-Null
$d.Content.Find.Execute("63+23=", $true, $false, $false, $false, $false, $true, 1, $false, "@@23@@", 2) | Out-Null
$d.Content.Find.Execute("96-3=", $true, $false, $false, $false, $false, $true, 1, $false, "@@24@@", 2) | Out-Null
$d.Content.Find.Execute("19-13=", $true, $false, $false, $false, $false, $true, 1, $false, "@@25@@", 2) | Out-Null
$d.Content.Find.Execute("71+20=", $true, $false, $false, $false, $false, $true, 1, $false, "@@26@@", 2) | Out-Null
$d.Content.Find.Execute("25+55=", $true, $false, $false, $false, $false, $true, 1, $false, "@@27@@", 2) | Out-Null
$d.Content.Find.Execute("68-17=", $true, $false, $false, $false, $false, $true, 1, $false, "@@28@@", 2) | Out-Null
$d.Content.Find.Execute("48-26=", $true, $false, $false, $false, $false, $true, 1, $false, "@@29@@", 2) | Out-Null
$d.Content.Find.Execute("44+54=", $true, $false, $false, $false, $false, $true, 1, $false, "@@30@@", 2) | Out-Null
$d.Content.Find.Execute("42+10=", $true, $false, $false, $false, $false, $true, 1, $false, "@@31@@", 2) | Out-Null
$d.Content.Find.Execute("37+35=", $true, $false, $false, $false, $false, $true, 1, $false, "@@32@@", 2) | Out-Null
$d.Content.Find.Execute("73-68=", $true, $false, $false, $false, $false, $true, 1, $false, "@@33@@", 2) | Out-Null
$d.Content.Find.Execute("74-52=", $true, $false, $false, $false, $false, $true, 1, $false, "@@34@@", 2) | Out-Null
$d.Content.Find.Execute("14+80=", $true, $false, $false, $false, $false, $true, 1, $false, "@@35@@", 2) | Out-Null
$d.Content.Find.Execute("30-13=", $true, $false, $false, $false, $false, $true, 1, $false, "@@36@@", 2) | Out-Null
$d.Content.Find.Execute("6+70=", $true, $false, $false, $false, $false, $true, 1, $false, "@@37@@", 2) | Out-Null
$d.Content.Find.Execute("44-37=", $true, $false, $false, $false, $false, $true, 1, $false, "@@38@@", 2) | Out-Null
$d.Content.Find.Execute("19+36=", $true, $false, $false, $false, $false, $true, 1, $false, "@@39@@", 2) | Out-Null
$d.Content.Find.Execute("9+29=", $true, $false, $false, $false, $false, $true, 1, $false, "@@40@@", 2) | Out-Null
$d.Content.Find.Execute("93-43=", $true, $false, $false, $false, $false, $true, 1, $false, "@@41@@", 2) | Out-Null
$d.Content.Find.Execute("22+41=", $true, $false, $false, $false, $false, $true, 1, $false, "@@42@@", 2) | Out-Null
$d.Content.Find.Execute("65-25=", $true, $false, $false, $false, $false, $true, 1, $false, "@@43@@", 2) | Out-Null
$d.Content.Find.Execute("58-39=", $true, $false, $false, $false, $false, $true, 1, $false, "@@44@@", 2) | Out-Null
$d.Content.Find.Execute("12+43=", $true, $false, $false, $false, $false, $true, 1, $false, "@@45@@", 2) | Out-Null
$d.Content.Find.Execute("99-69=", $true, $false, $false, $false, $false, $true, 1, $false, "@@46@@", 2) | Out-Null
$d.Content.Find.Execute("26-3=", $true, $false, $false, $false, $false, $true, 1, $false, "@@47@@", 2) | Out-Null
$d.Content.Find.Execute("83-50=", $true, $false, $false, $false, $false, $true, 1, $false, "@@48@@", 2) | Out-Null
$d.Content.Find.Execute("35-16=", $true, $false, $false, $false, $false, $true, 1, $false, "@@49@@", 2) | Out-Null
$d.Content.Find.Execute("78-35=", $true, $false, $false, $false, $false, $true, 1, $false, "@@50@@", 2) | Out-Null
$d.Content.Find.Execute("4+41=", $true, $false, $false, $false, $false, $true, 1, $false, "@@51@@", 2) | Out-Null
$d.Content.Find.Execute("25+11=", $true, $false, $false, $false, $false, $true, 1, $false, "@@52@@", 2) | Out-Null
$d.Content.Find.Execute("58+21=", $true, $false, $false, $false, $false, $true, 1, $false, "@@53@@", 2) | Out-Null
$d.Content.Find.Execute("72+18=", $true, $false, $false, $false, $false, $true, 1, $false, "@@54@@", 2) | Out-Null
$d.Content.Find.Execute("84-3=", $true, $false, $false, $false, $false, $true, 1, $false, "@@55@@", 2) | Out-Null
$d.Content.Find.Execute("11+38=", $true, $false, $false, $false, $false, $true, 1, $false, "@@56@@", 2) | Out-Null
$d.Content.Find.Execute("86-72=", $true, $false, $false, $false, $false, $true, 1, $false, "@@57@@", 2) | Out-Null
$d.Content.Find.Execute("65-31=", $true, $false, $false, $false, $false, $true, 1, $false, "@@58@@", 2) | Out-Null
$d.Content.Find.Execute("75-3=", $true, $false, $false, $false, $false, $true, 1, $false, "@@59@@", 2) | Out-Null
$d.Content.Find.Execute("96-76=", $true, $false, $false, $false, $false, $true, 1, $false, "@@60@@", 2) | Out-Null
$d.Content.Find.Execute("22+76=", $true, $false, $false, $false, $false, $true, 1, $false, "@@61@@", 2) | Out-Null
$d.Content.Find.Execute("26+39=", $true, $false, $false, $false, $false, $true, 1, $false, "@@62@@", 2) | Out-Null
$d.Content.Find.Execute("12+5=", $true, $false, $false, $false, $false, $true, 1, $false, "@@63@@", 2) | Out-Null
$d.Content.Find.Execute("68+15=", $true, $false, $false, $false, $false, $true, 1, $false, "@@64@@", 2) | Out-Null
$d.Content.Find.Execute("34-17=", $true, $false, $false, $false, $false, $true, 1, $false, "@@65@@", 2) | Out-Null
$d.Content.Find.Execute("80-77=", $true, $false, $false, $false, $false, $true, 1, $false, "@@66@@", 2) | Out-Null
$d.Content.Find.Execute("0+92=", $true, $false, $false, $false, $false, $true, 1, $false, "@@67@@", 2) | Out-Null
$d.Content.Find.Execute("94-34=", $true, $false, $false, $false, $false, $true, 1, $false, "@@68@@", 2) | Out-Null
$d.Content.Find.Execute("3+18=", $true, $false, $false, $false, $false, $true, 1, $false, "@@69@@", 2) | Out-Null
$d.Content.Find.Execute("76-58=", $true, $false, $false, $false, $false, $true, 1, $false, "@@70@@", 2) | Out-Null
$d.Content.Find.Execute("17+13=", $true, $false, $false, $false, $false, $true, 1, $false, "@@71@@", 2) | Out-Null
$d.Content.Find.Execute("51-50=", $true, $false, $false, $false, $false, $true, 1, $false, "@@72@@", 2) | Out-Null
$d.Content.Find.Execute("9+10=", $true, $false, $false, $false, $false, $true, 1, $false, "@@73@@", 2) | Out-Null
$d.Content.Find.Execute("58+9=", $true, $false, $false, $false, $false, $true, 1, $false, "@@74@@", 2) | Out-Null
$d.Content.Find.Execute("20+63=", $true, $false, $false, $false, $false, $true, 1, $false, "@@75@@", 2) | Out-Null
$d.Content.Find.Execute("34+5=", $true, $false, $false, $false, $false, $true, 1, $false, "@@76@@", 2) | Out-Null
$d.Content.Find.Execute("1-0=", $true, $false, $false, $false, $false, $true, 1, $false, "@@77@@", 2) | Out-Null
$d.Content.Find.Execute("6+75=", $true, $false, $false, $false, $false, $true, 1, $false, "@@78@@", 2) | Out-Null
$d.Content.Find.Execute("64-47=", $true, $false, $false, $false, $false, $true, 1, $false, "@@79@@", 2) | Out-Null
$d.Content.Find.Execute("44-35=", $true, $false, $false, $false, $false, $true, 1, $false, "@@80@@", 2) | Out-Null
$d.Content.Find.Execute("79-12=", $true, $false, $false, $false, $false, $true, 1, $false, "@@81@@", 2) | Out-Null
$d.Content.Find.Execute("97-85=", $true, $false, $false, $false, $false, $true, 1, $false, "@@82@@", 2) | Out-Null
$d.Content.Find.Execute("88-82=", $true, $false, $false, $false, $false, $true, 1, $false, "@@83@@", 2) | Out-Null
$d.Content.Find.Execute("4+48=", $true, $false, $false, $false, $false, $true, 1, $false, "@@84@@", 2) | Out-Null
$d.Content.Find.Execute("20-16=", $true, $false, $false, $false, $false, $true, 1, $false, "@@85@@", 2) | Out-Null
$d.Content.Find.Execute("37+51=", $true, $false, $false, $false, $false, $true, 1, $false, "@@86@@", 2) | Out-Null
$d.Content.Find.Execute("5+71=", $true, $false, $false, $false, $false, $true, 1, $false, "@@87@@", 2) | Out-Null
$d.Content.Find.Execute("94-43=", $true, $false, $false, $false, $false, $true, 1, $false, "@@88@@", 2) | Out-Null
$d.Content.Find.Execute("56-14=", $true, $false, $false, $false, $false, $true, 1, $false, "@@89@@", 2) | Out-Null
$d.Content.Find.Execute("25+70=", $true, $false, $false, $false, $false, $true, 1, $false, "@@90@@", 2) | Out-Null
$d.Content.Find.Execute("9+85=", $true, $false, $false, $false, $false, $true, 1, $false, "@@91@@", 2) | Out-Null
$d.Content.Find.Execute("4+24=", $true, $false, $false, $false, $false, $true, 1, $false, "@@92@@", 2) | Out-Null
$d.Content.Find.Execute("91-47=", $true, $false, $false, $false, $false, $true, 1, $false, "@@93@@", 2) | Out-Null
$d.Content.Find.Execute("92-91=", $true, $false, $false, $false, $false, $true, 1, $false, "@@94@@", 2) | Out-Null
$d.Content.Find.Execute("69-11=", $true, $false, $false, $false, $false, $true, 1, $false, "@@95@@", 2) | Out-Null
$d.Content.Find.Execute("47-4=", $true, $false, $false, $false, $false, $true, 1, $false, "@@96@@", 2) | Out-Null
$d.Content.Find.Execute("42+20=", $true, $false, $false, $false, $false, $true, 1, $false, "@@97@@", 2) | Out-Null
$d.Content.Find.Execute("42+25=", $true, $false, $false, $false, $false, $true, 1, $false, "@@98@@", 2) | Out-Null
$d.Content.Find.Execute("34+22=", $true, $false, $false, $false, $false, $true, 1, $false, "@@99@@", 2) | Out-Null

# Phase 2: replace each placeholder token with its final new value
$d.Content.Find.Execute("@@0@@", $true, $false, $false, $false, $false, $true, 1, $false, "34-20=", 2) | Out-Null
$d.Content.Find.Execute("@@1@@", $true, $false, $false, $false, $false, $true, 1, $false, "46-29=", 2) | Out-Null
$d.Content.Find.Execute("@@2@@", $true, $false, $false, $false, $false, $true, 1, $false, "24+64=", 2) | Out-Null
$d.Content.Find.Execute("@@3@@", $true, $false, $false, $false, $false, $true, 1, $false, "29+27=", 2) | Out-Null
$d.Content.Find.Execute("@@4@@", $true, $false, $false, $false, $false, $true, 1, $false, "24+14=", 2) | Out-Null
$d.Content.Find.Execute("@@5@@", $true, $false, $false, $false, $false, $true, 1, $false, "17+63=", 2) | Out-Null
$d.Content.Find.Execute("@@6@@", $true, $false, $false, $false, $false, $true, 1, $false, "83-29=", 2) | Out-Null
$d.Content.Find.Execute("@@7@@", $true, $false, $false, $false, $false, $true, 1, $false, "74+17=", 2) | Out-Null
$d.Content.Find.Execute("@@8@@", $true, $false, $false, $false, $false, $true, 1, $false, "45+24=", 2) | Out-Null
$d.Content.Find.Execute("@@9@@", $true, $false, $false, $false, $false, $true, 1, $false, "75-22=", 2) | Out-Null
$d.Content.Find.Execute("@@10@@", $true, $false, $false, $false, $false, $true, 1, $false, "30+57=", 2) | Out-Null
$d.Content.Find.Execute("@@11@@", $true, $false, $false, $false, $false, $true, 1, $false, "66-53=", 2) | Out-Null
$d.Content.Find.Execute("@@12@@", $true, $false, $false, $false, $false, $true, 1, $false, "71-41=", 2) | Out-Null
$d.Content.Find.Execute("@@13@@", $true, $false, $false, $false, $false, $true, 1, $false, "28+58=", 2) | Out-Null
$d.Content.Find.Execute("@@14@@", $true, $false, $false, $false, $false, $true, 1, $false, "44+15=", 2) | Out-Null
$d.Content.Find.Execute("@@15@@", $true, $false, $false, $false, $false, $true, 1, $false, "63+7=", 2) | Out-Null
$d.Content.Find.Execute("@@16@@", $true, $false, $false, $false, $false, $true, 1, $false, "57+15=", 2) | Out-Null
$d.Content.Find.Execute("@@17@@", $true, $false, $false, $false, $false, $true, 1, $false, "86+1=", 2) | Out-Null
$d.Content.Find.Execute("@@18@@", $true, $false, $false, $false, $false, $true, 1, $false, "77-17=", 2) | Out-Null
$d.Content.Find.Execute("@@19@@", $true, $false, $false, $false, $false, $true, 1, $false, "41+55=", 2) | Out-Null
$d.Content.Find.Execute("@@20@@", $true, $false, $false, $false, $false, $true, 1, $false, "45-11=", 2) | Out-Null
$d.Content.Find.Execute("@@21@@", $true, $false, $false, $false, $false, $true, 1, $false, "23-20=", 2) | Out-Null
$d.Content.Find.Execute("@@22@@", $true, $false, $false, $false, $false, $true, 1, $false, "47+44=", 2) | Out-Null
$d.Content.Find.Execute("@@23@@", $true, $false, $false, $false, $false, $true, 1, $false, "0+64=", 2) | Out-Null
$d.Content.Find.Execute("@@24@@", $true, $false, $false, $false, $false, $true, 1, $false, "64+2=", 2) | Out-Null
$d.Content.Find.Execute("@@25@@", $true, $false, $false, $false, $false, $true, 1, $false, "57+20=", 2) | Out-Null
$d.Content.Find.Execute("@@26@@", $true, $false, $false, $false, $false, $true, 1, $false, "83-55=", 2) | Out-Null
$d.Content.Find.Execute("@@27@@", $true, $false, $false, $false, $false, $true, 1, $false, "98-57=", 2) | Out-Null
$d.Content.Find.Execute("@@28@@", $true, $false, $false, $false, $false, $true, 1, $false, "40-3=", 2) | Out-Null
$d.Content.Find.Execute("@@29@@", $true, $false, $false, $false, $false, $true, 1, $false, "95-89=", 2) | Out-Null
$d.Content.Find.Execute("@@30@@", $true, $false, $false, $false, $false, $true, 1, $false, "78-8=", 2) | Out-Null
$d.Content.Find.Execute("@@31@@", $true, $false, $false, $false, $false, $true, 1, $false, "16-11=", 2) | Out-Null
$d.Content.Find.Execute("@@32@@", $true, $false, $false, $false, $false, $true, 1, $false, "83-0=", 2) | Out-Null
$d.Content.Find.Execute("@@33@@", $true, $false, $false, $false, $false, $true, 1, $false, "68-56=", 2) | Out-Null
$d.Content.Find.Execute("@@34@@", $true, $false, $false, $false, $false, $true, 1, $false, "29+42=", 2) | Out-Null
$d.Content.Find.Execute("@@35@@", $true, $false, $false, $false, $false, $true, 1, $false, "71+21=", 2) | Out-Null
$d.Content.Find.Execute("@@36@@", $true, $false, $false, $false, $false, $true, 1, $false, "73-3=", 2) | Out-Null
$d.Content.Find.Execute("@@37@@", $true, $false, $false, $false, $false, $true, 1, $false, "39-35=", 2) | Out-Null
$d.Content.Find.Execute("@@38@@", $true, $false, $false, $false, $false, $true, 1, $false, "74+9=", 2) | Out-Null
$d.Content.Find.Execute("@@39@@", $true, $false, $false, $false, $false, $true, 1, $false, "34-0=", 2) | Out-Null
$d.Content.Find.Execute("@@40@@", $true, $false, $false, $false, $false, $true, 1, $false, "2+23=", 2) | Out-Null
$d.Content.Find.Execute("@@41@@", $true, $false, $false, $false, $false, $true, 1, $false, "81-59=", 2) | Out-Null
$d.Content.Find.Execute("@@42@@", $true, $false, $false, $false, $false, $true, 1, $false, "22+65=", 2) | Out-Null
$d.Content.Find.Execute("@@43@@", $true, $false, $false, $false, $false, $true, 1, $false, "21+32=", 2) | Out-Null
$d.Content.Find.Execute("@@44@@", $true, $false, $false, $false, $false, $true, 1, $false, "83-58=", 2) | Out-Null
$d.Content.Find.Execute("@@45@@", $true, $false, $false, $false, $false, $true, 1, $false, "74-31=", 2) | Out-Null
$d.Content.Find.Execute("@@46@@", $true, $false, $false, $false, $false, $true, 1, $false, "56+29=", 2) | Out-Null
$d.Content.Find.Execute("@@47@@", $true, $false, $false, $false, $false, $true, 1, $false, "45-4=", 2) | Out-Null
$d.Content.Find.Execute("@@48@@", $true, $false, $false, $false, $false, $true, 1, $false, "52-51=", 2) | Out-Null
$d.Content.Find.Execute("@@49@@", $true, $false, $false, $false, $false, $true, 1, $false, "0+87=", 2) | Out-Null
$d.Content.Find.Execute("@@50@@", $true, $false, $false, $false, $false, $true, 1, $false, "84-1=", 2) | Out-Null
$d.Content.Find.Execute("@@51@@", $true, $false, $false, $false, $false, $true, 1, $false, "47+34=", 2) | Out-Null
$d.Content.Find.Execute("@@52@@", $true, $false, $false, $false, $false, $true, 1, $false, "58-19=", 2) | Out-Null
$d.Content.Find.Execute("@@53@@", $true, $false, $false, $false, $false, $true, 1, $false, "28-8=", 2) | Out-Null
$d.Content.Find.Execute("@@54@@", $true, $false, $false, $false, $false, $true, 1, $false, "17+46=", 2) | Out-Null
$d.Content.Find.Execute("@@55@@", $true, $false, $false, $false, $false, $true, 1, $false, "97-1=", 2) | Out-Null
$d.Content.Find.Execute("@@56@@", $true, $false, $false, $false, $false, $true, 1, $false, "67-18=", 2) | Out-Null
$d.Content.Find.Execute("@@57@@", $true, $false, $false, $false, $false, $true, 1, $false, "4+78=", 2) | Out-Null
$d.Content.Find.Execute("@@58@@", $true, $false, $false, $false, $false, $true, 1, $false, "63+6=", 2) | Out-Null
$d.Content.Find.Execute("@@59@@", $true, $false, $false, $false, $false, $true, 1, $false, "16+13=", 2) | Out-Null
$d.Content.Find.Execute("@@60@@", $true, $false, $false, $false, $false, $true, 1, $false, "82-33=", 2) | Out-Null
$d.Content.Find.Execute("@@61@@", $true, $false, $false, $false, $false, $true, 1, $false, "67+9=", 2) | Out-Null
$d.Content.Find.Execute("@@62@@", $true, $false, $false, $false, $false, $true, 1, $false, "25+47=", 2) | Out-Null
$d.Content.Find.Execute("@@63@@", $true, $false, $false, $false, $false, $true, 1, $false, "69-34=", 2) | Out-Null
$d.Content.Find.Execute("@@64@@", $true, $false, $false, $false, $false, $true, 1, $false, "49+29=", 2) | Out-Null
$d.Content.Find.Execute("@@65@@", $true, $false, $false, $false, $false, $true, 1, $false, "27+47=", 2) | Out-Null
$d.Content.Find.Execute("@@66@@", $true, $false, $false, $false, $false, $true, 1, $false, "71-25=", 2) | Out-Null
$d.Content.Find.Execute("@@67@@", $true, $false, $false, $false, $false, $true, 1, $false, "69-15=", 2) | Out-Null
$d.Content.Find.Execute("@@68@@", $true, $false, $false, $false, $false, $true, 1, $false, "60-31=", 2) | Out-Null
$d.Content.Find.Execute("@@69@@", $true, $false, $false, $false, $false, $true, 1, $false, "17+77=", 2) | Out-Null
$d.Content.Find.Execute("@@70@@", $true, $false, $false, $false, $false, $true, 1, $false, "17+35=", 2) | Out-Null
$d.Content.Find.Execute("@@71@@", $true, $false, $false, $false, $false, $true, 1, $false, "81-76=", 2) | Out-Null
$d.Content.Find.Execute("@@72@@", $true, $false, $false, $false, $false, $true, 1, $false, "19+6=", 2) | Out-Null
$d.Content.Find.Execute("@@73@@", $true, $false, $false, $false, $false, $true, 1, $false, "97-11=", 2) | Out-Null
$d.Content.Find.Execute("@@74@@", $true, $false, $false, $false, $false, $true, 1, $false, "79-47=", 2) | Out-Null
$d.Content.Find.Execute("@@75@@", $true, $false, $false, $false, $false, $true, 1, $false, "26+28=", 2) | Out-Null
$d.Content.Find.Execute("@@76@@", $true, $false, $false, $false, $false, $true, 1, $false, "95-58=", 2) | Out-Null
$d.Content.Find.Execute("@@77@@", $true, $false, $false, $false, $false, $true, 1, $false, "91-2=", 2) | Out-Null
$d.Content.Find.Execute("@@78@@", $true, $false, $false, $false, $false, $true, 1, $false, "48+3=", 2) | Out-Null
$d.Content.Find.Execute("@@79@@", $true, $false, $false, $false, $false, $true, 1, $false, "61-16=", 2) | Out-Null
$d.Content.Find.Execute("@@80@@", $true, $false, $false, $false, $false, $true, 1, $false, "90+1=", 2) | Out-Null
$d.Content.Find.Execute("@@81@@", $true, $false, $false, $false, $false, $true, 1, $false, "14+43=", 2) | Out-Null
$d.Content.Find.Execute("@@82@@", $true, $false, $false, $false, $false, $true, 1, $false, "29-16=", 2) | Out-Null
$d.Content.Find.Execute("@@83@@", $true, $false, $false, $false, $false, $true, 1, $false, "79-5=", 2) | Out-Null
$d.Content.Find.Execute("@@84@@", $true, $false, $false, $false, $false, $true, 1, $false, "25+70=", 2) | Out-Null
$d.Content.Find.Execute("@@85@@", $true, $false, $false, $false, $false, $true, 1, $false, "60-23=", 2) | Out-Null
$d.Content.Find.Execute("@@86@@", $true, $false, $false, $false, $false, $true, 1, $false, "61+27=", 2) | Out-Null
$d.Content.Find.Execute("@@87@@", $true, $false, $false, $false, $false, $true, 1, $false, "0+41=", 2) | Out-Null
$d.Content.Find.Execute("@@88@@", $true, $false, $false, $false, $false, $true, 1, $false, "90-43=", 2) | Out-Null
$d.Content.Find.Execute("@@89@@", $true, $false, $false, $false, $false, $true, 1, $false, "33-2=", 2) | Out-Null
$d.Content.Find.Execute("@@90@@", $true, $false, $false, $false, $false, $true, 1, $false, "93-59=", 2) | Out-Null
$d.Content.Find.Execute("@@91@@", $true, $false, $false, $false, $false, $true, 1, $false, "35-9=", 2) | Out-Null
$d.Content.Find.Execute("@@92@@", $true, $false, $false, $false, $false, $true, 1, $false, "96-25=", 2) | Out-Null
$d.Content.Find.Execute("@@93@@", $true, $false, $false, $false, $false, $true, 1, $false, "44+3=", 2) | Out-Null
$d.Content.Find.Execute("@@94@@", $true, $false, $false, $false, $false, $true, 1, $false, "70+9=", 2) | Out-Null
$d.Content.Find.Execute("@@95@@", $true, $false, $false, $false, $false, $true, 1, $false, "12+27=", 2) | Out-Null
$d.Content.Find.Execute("@@96@@", $true, $false, $false, $false, $false, $true, 1, $false, "5+85=", 2) | Out-Null
$d.Content.Find.Execute("@@97@@", $true, $false, $false, $false, $false, $true, 1, $false, "23+71=", 2) | Out-Null
$d.Content.Find.Execute("@@98@@", $true, $false, $false, $false, $false, $true, 1, $false, "78-14=", 2) | Out-Null
$d.Content.Find.Execute("@@99@@", $true, $false, $false, $false, $false, $true, 1, $false, "3+36=", 2) | Out-Null
